# The workbook's single sheet was renamed from "O4_sps" to "Q4_sps"
# (the workbook/file itself is "Q4_sps.xlsx", so the tab name was
# corrected to match), and the live cell selection/cursor position was
# moved from G23:G24 (active cell G23) to the single cell L26.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Q4_sps"
$ws.Range("L26").Select()
